# Update cryptocurrency price/volume figures (and restore the
# ARBITRUM / ImmutableX row order) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces text (matches the source cells, which are
# stored as inline strings, not numbers) so values like "216.09"
# or "27.122.28" are not auto-converted to numeric cells.
$ws.Range("D2").Value = "'27.122.28"
$ws.Range("E2").Value = "'  -0.23%  "
$ws.Range("D3").Value = "'1.628.94"
$ws.Range("E3").Value = "'  -1.10%  "
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("D5").Value = "'216.09"
$ws.Range("E5").Value = "'  -1.12%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("E6").Value = "'  +0.89%  "
$ws.Range("E7").Value = "'  -0.10%  "
$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "'  -1.43%  "
$ws.Range("D9").Value = "'0.0623"
$ws.Range("E9").Value = "'  -0.86%  "
$ws.Range("D10").Value = "'20.04"
$ws.Range("E10").Value = "'  +0.01%  "
$ws.Range("E11").Value = "'  +0.35%  "
$ws.Range("D12").Value = "'1.857.18"
$ws.Range("E12").Value = "'  -1.10%  "
$ws.Range("D13").Value = "'1.641.54"
$ws.Range("E13").Value = "'  -0.33%  "
$ws.Range("E14").Value = "'  -0.64%  "
$ws.Range("D15").Value = "'0.538"
$ws.Range("E15").Value = "'  +0.13%  "
$ws.Range("D16").Value = "'64.95"
$ws.Range("E16").Value = "'  -3.62%  "
$ws.Range("D17").Value = "'27.080.74"
$ws.Range("D18").Value = "'0.0₃0731"
$ws.Range("E18").Value = "'  -1.31%  "
$ws.Range("D19").Value = "'213.55"
$ws.Range("E19").Value = "'  -2.74%  "
$ws.Range("E20").Value = "'  +0.01%  "
$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = "'  +0.54%  "
$ws.Range("D22").Value = "'4.37"
$ws.Range("E22").Value = "'  -1.59%  "
$ws.Range("D23").Value = "'2.46"
$ws.Range("E23").Value = "'  -1.53%  "
$ws.Range("D24").Value = "'9.05"
$ws.Range("E24").Value = "'  -1.84%  "
$ws.Range("D25").Value = "'147.42"
$ws.Range("E25").Value = "'  -0.61%  "
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("D27").Value = "'7.31"
$ws.Range("E27").Value = "'  -1.35%  "
$ws.Range("E28").Value = "'  -0.93%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "'  -1.68%  "
$ws.Range("E30").Value = "'  -0.43%  "
$ws.Range("E31").Value = "'  -1.18%  "
$ws.Range("D32").Value = "'3.37"
$ws.Range("E32").Value = "'  +0.36%  "
$ws.Range("D33").Value = "'2.99"
$ws.Range("E33").Value = "'  -1.60%  "
$ws.Range("D34").Value = "'1.312.33"
$ws.Range("E34").Value = "'  +3.37%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "'  -1.83%  "
$ws.Range("E36").Value = "'  -0.42%  "
$ws.Range("E37").Value = "'  -1.98%  "
$ws.Range("B38").Value = "'ARBITRUM"
$ws.Range("C38").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'0.841"
$ws.Range("E38").Value = "'  -0.55%  "
$ws.Range("B39").Value = "'ImmutableX"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.536"
$ws.Range("E39").Value = "'  -1.24%  "
$ws.Range("E40").Value = "'  -0.10%  "
$ws.Range("D41").Value = "'2.27"
$ws.Range("E41").Value = "'  +1.99%  "
$ws.Range("D42").Value = "'0.802"
$ws.Range("E42").Value = "'  -1.03%  "
$ws.Range("D43").Value = "'5.24"
$ws.Range("E43").Value = "'  -2.53%  "
$ws.Range("D44").Value = "'1.766.21"
$ws.Range("E44").Value = "'  -1.26%  "
$ws.Range("D45").Value = "'62.56"
$ws.Range("E45").Value = "'  +0.08%  "
$ws.Range("D46").Value = "'90.69"
$ws.Range("E46").Value = "'  -1.82%  "
$ws.Range("E47").Value = "'  -0.25%  "
$ws.Range("E48").Value = "'  +18.13%  "
$ws.Range("D49").Value = "'0.796"
$ws.Range("E49").Value = "'  +17.62%  "
$ws.Range("E50").Value = "'  -0.07%  "
$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = "'  -1.82%  "
